$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 23:22"

# Update Cataluña row (row 5) statistics
$ws.Range("B5").Value = 34027
$ws.Range("C5").Value = 14975
$ws.Range("D5").Value = 15610
$ws.Range("E5").Value = 3442
